$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Yes" marker cell in C10 (row 10), copying the style used by
# existing "Yes" cells (e.g. C4) so the green fill / font carries over.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C10").Value = "Yes"

# Row 22: add "Yes" marker + completion date note, and bump the row height
# like the other "REQUIRES POLISHING" rows (e.g. row 21).
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = "Yes"

$ws.Range("D21").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null

$d22 = $ws.Range("D22")
$d22.Value = "02/06/2024 (REQUIRES POLISHING)"
$d22.Characters(13, 18).Font.Bold = $true

$ws.Rows.Item(22).RowHeight = 15

# Row 27: same treatment as row 22, reusing D22's rich text so both cells
# point at the same shared string entry.
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = "Yes"

$ws.Range("D22").Copy() | Out-Null
$ws.Paste($ws.Range("D27")) | Out-Null

$ws.Rows.Item(27).RowHeight = 15

# Update selection to match the saved cursor position in the source file.
$ws.Range("E30").Select() | Out-Null
